# Estonia Meistriliiga.xlsx - "Atualizacao de bases das ligas" update
#
# This update does two things:
#  1. Corrects the team names that had been swapped in the lookup/shared
#     string table: "JK Tammeka Tartu" and "JK Tallinna Kalev" had been
#     mixed up, so the two names are swapped wherever they occur.
#  2. Because of (1), three pairs of rows that hold match data for these
#     two teams need their row content (everything except the HomeTeam
#     column, which is handled by the name swap itself / or which simply
#     moves with the row) swapped between the two rows so each match's
#     statistics line up with the correct fixture again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param(
        [int]$Row1,
        [int]$Row2,
        [int[]]$Columns
    )

    # Snapshot both rows first so we don't clobber values we still need to read.
    $row1Values = @{}
    $row2Values = @{}
    foreach ($col in $Columns) {
        $row1Values[$col] = $ws.Cells.Item($Row1, $col).Value2
        $row2Values[$col] = $ws.Cells.Item($Row2, $col).Value2
    }
    foreach ($col in $Columns) {
        $ws.Cells.Item($Row1, $col).Value = $row2Values[$col]
        $ws.Cells.Item($Row2, $col).Value = $row1Values[$col]
    }
}

# Columns B (id) and G..AC (everything after AwayTeam) swap between rows 4
# and 5; column F (HomeTeam) is left untouched on this pair.
$colsNoHomeTeam = 2,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29
Swap-RowRange -Row1 4 -Row2 5 -Columns $colsNoHomeTeam

# Columns B and F..AC (including HomeTeam) swap between rows 95 and 96, and
# between rows 105 and 107.
$colsWithHomeTeam = 2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29
Swap-RowRange -Row1 95 -Row2 96 -Columns $colsWithHomeTeam
Swap-RowRange -Row1 105 -Row2 107 -Columns $colsWithHomeTeam
